$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F (Person Responsible), shifting old F (Date Last
# Edited) to G and old G (Comments) to H.
$ws.Columns("F:F").Insert()

# Update the "Format" version string.
$ws.Range("C2").Value = "v0.2.0"

# Populate the newly inserted "Person Responsible" column.
$ws.Range("F3").Value = "Person Responsible"
$ws.Range("F4").Value = "personResponsible"
$ws.Range("F5").Value = "Person responsible that added this source and the corresponding entries"
$ws.Range("F6").Value = "-"
$ws.Range("F7").Value = "[text]"
$ws.Range("F8").Value = "Tester"
$ws.Range("F9").Value = "Tester"

# Fix up the column widths: column F keeps original A/C widths slightly
# adjusted, and the new column F should match the old F (now G) width.
$ws.Columns("A:A").ColumnWidth = 5.5703125
$ws.Columns("C:C").ColumnWidth = 24.5703125
$ws.Columns("F:F").ColumnWidth = 25.7109375

Write-Host "done"
